$wb = $excel.ActiveWorkbook

# --- Restructure "MahanaAamdanDepositInputter": insert 4 new columns (C:F) ---
# This shifts the existing T.C (Azure)/T.C (Desc.)/Error columns from C,D,E to G,H,I.
$ws = $wb.Worksheets.Item("MahanaAamdanDepositInputter")
[void]$ws.Range("C1:F1").EntireColumn.Insert()

# --- Update the data row ---
# "1M" is written first so its shared-string slot is allocated before the
# Auth sheet's string and the new header strings (matches authoring order).
$ws = $wb.Worksheets.Item("MahanaAamdanDepositInputter")
$ws.Range("B2").Value = "1M"
$ws.Range("A2").Value = 16988322

# --- Add the new "MahanaAamdanDepositAuth" sheet right after the Inputter sheet ---
$ws = $wb.Worksheets.Item("MahanaAamdanDepositInputter")
$authSheet = $wb.Worksheets.Add($null, $ws)
$authSheet.Name = "MahanaAamdanDepositAuth"
$authSheet.Range("A1").Value = "value:1:1:1"
$authSheet.Range("A2").Value = 1007887136

# --- Fill in the new header cells on the Inputter sheet ---
$ws = $wb.Worksheets.Item("MahanaAamdanDepositInputter")
$ws.Range("C1").Value = "CUST.REMARKS:1"
$ws.Range("D1").Value = "FIQAH"
$ws.Range("E1").Value = "INTEND.DATE"
$ws.Range("F1").Value = "EXP.DATE"

# --- Add the new blank "Sheet1" before the Inputter sheet ---
$ws = $wb.Worksheets.Item("MahanaAamdanDepositInputter")
$sheet1 = $wb.Worksheets.Add($ws)

# --- Column widths on the Inputter sheet (closest reproducible approximation) ---
$ws = $wb.Worksheets.Item("MahanaAamdanDepositInputter")
$ws.Columns.Item(1).ColumnWidth = 11.5
$ws.Columns.Item(2).ColumnWidth = 10.333333333333334
$ws.Columns.Item(3).ColumnWidth = 20.666666666666668
$ws.Columns.Item(4).ColumnWidth = 7.666666666666667
$ws.Columns.Item(5).ColumnWidth = 16
$ws.Columns.Item(6).ColumnWidth = 11.333333333333334
$ws.Columns.Item(7).ColumnWidth = 13
$ws.Columns.Item(8).ColumnWidth = 12.333333333333334
$ws.Columns.Item(9).ColumnWidth = 6.166666666666667

# --- Column width on the Auth sheet ---
$authSheet = $wb.Worksheets.Item("MahanaAamdanDepositAuth")
$authSheet.Columns.Item(1).ColumnWidth = 10.166666666666666

# --- Selections / active sheet ---
$authSheet = $wb.Worksheets.Item("MahanaAamdanDepositAuth")
[void]$authSheet.Select()
[void]$authSheet.Range("B1").Select()

$ws = $wb.Worksheets.Item("MahanaAamdanDepositInputter")
[void]$ws.Select()
[void]$ws.Range("C6").Select()

Write-Output "done"
